$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1135171.8
$ws.Range("J17").Value = 1689455
$ws.Range("L17").Value = 5068365
$ws.Range("N17").Value = -5068701
$ws.Range("H70").Value = 999
$ws.Range("J70").Value = 999
$ws.Range("L70").Value = 2997
$ws.Range("N70").Value = -3537
$ws.Range("H73").Value = 999
$ws.Range("J73").Value = 999
$ws.Range("L73").Value = 2997
$ws.Range("N73").Value = -4869
$ws.Range("H112").Value = 6077.871
$ws.Range("I112").Value = 340
$ws.Range("J112").Value = 7181.3076
$ws.Range("K112").Value = 1020
$ws.Range("L112").Value = 21543.9228
$ws.Range("M112").Value = 88
$ws.Range("N112").Value = -23759.9228
$ws.Range("H137").Value = 2714.9607
$ws.Range("I137").Value = 2204.4055
$ws.Range("J137").Value = 4064.2856
$ws.Range("K137").Value = 6613.2165
$ws.Range("L137").Value = 12192.8568
$ws.Range("M137").Value = -4063.2165
$ws.Range("N137").Value = -17292.8568
$ws.Range("H138").Value = 2433.9565
$ws.Range("J138").Value = 3500.182
$ws.Range("L138").Value = 10500.546
$ws.Range("N138").Value = -20780.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2879.1177
$ws.Range("I61").Value = 2398.8333
$ws.Range("J61").Value = 4031.8
$ws.Range("K61").Value = 2398.8333
$ws.Range("L61").Value = 4031.8
$ws.Range("M61").Value = -2186.8333
$ws.Range("N61").Value = -4455.8
$ws.Range("H74").Value = 2199.5625
$ws.Range("I74").Value = 1293.1428
$ws.Range("K74").Value = 1293.1428
$ws.Range("M74").Value = -419.1428000000001
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51748
$ws.Range("H77").Value = 2199.5625
$ws.Range("I77").Value = 1293.1428
$ws.Range("K77").Value = 6465.714
$ws.Range("M77").Value = -2097.714
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -158736
$ws.Range("H136").Value = 2879.1177
$ws.Range("I136").Value = 2398.8333
$ws.Range("J136").Value = 4031.8
$ws.Range("K136").Value = 7196.499899999999
$ws.Range("L136").Value = 12095.4
$ws.Range("M136").Value = -4646.499899999999
$ws.Range("N136").Value = -17195.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 287.5
$ws.Range("I7").Value = 400
$ws.Range("J7").Value = 175
$ws.Range("K7").Value = 400
$ws.Range("L7").Value = 175
$ws.Range("M7").Value = -287
$ws.Range("N7").Value = -401
$ws.Range("H31").Value = 9708.968000000001
$ws.Range("I31").Value = 1601.8462
$ws.Range("J31").Value = 15564.111
$ws.Range("K31").Value = 1601.8462
$ws.Range("L31").Value = 15564.111
$ws.Range("M31").Value = -1306.8462
$ws.Range("N31").Value = -16154.111
$ws.Range("H34").Value = 9708.968000000001
$ws.Range("I34").Value = 1601.8462
$ws.Range("J34").Value = 15564.111
$ws.Range("K34").Value = 1601.8462
$ws.Range("L34").Value = 15564.111
$ws.Range("M34").Value = -1399.8462
$ws.Range("N34").Value = -15968.111
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376
$ws.Range("N62").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("N65").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H69").Value = 27091
$ws.Range("I69").Value = 27091
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 27091
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -26342
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 34800
$ws.Range("J70").Value = 34800
$ws.Range("L70").Value = 34800
$ws.Range("N70").Value = -35430
$ws.Range("H72").Value = 27091
$ws.Range("I72").Value = 27091
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 81273
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -77529
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 34800
$ws.Range("J73").Value = 34800
$ws.Range("L73").Value = 34800
$ws.Range("N73").Value = -36984
$ws.Range("H88").Value = 17999
$ws.Range("J88").Value = 17999
$ws.Range("L88").Value = 17999
$ws.Range("N88").Value = -18811
$ws.Range("H91").Value = 17999
$ws.Range("J91").Value = 17999
$ws.Range("L91").Value = 17999
$ws.Range("N91").Value = -20807
$ws.Range("H132").Value = 4506388.5
$ws.Range("I132").Value = 1352.2
$ws.Range("K132").Value = 4056.6
$ws.Range("M132").Value = -1526.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 802.15
$ws.Range("J5").Value = 1243.2858
$ws.Range("L5").Value = 3729.8574
$ws.Range("N5").Value = -3953.8574
$ws.Range("H33").Value = 7226.7856
$ws.Range("I33").Value = 9192.272000000001
$ws.Range("J33").Value = 20
$ws.Range("K33").Value = 55153.63200000001
$ws.Range("L33").Value = 120
$ws.Range("M33").Value = -54870.63200000001
$ws.Range("N33").Value = -686
$ws.Range("H82").Value = 1000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H88").Value = 11346.833
$ws.Range("J88").Value = 11346.833
$ws.Range("L88").Value = 34040.499
$ws.Range("N88").Value = -34896.499
$ws.Range("H91").Value = 11346.833
$ws.Range("J91").Value = 11346.833
$ws.Range("L91").Value = 34040.499
$ws.Range("N91").Value = -37004.499
$ws.Range("H104").Value = 4125.8
$ws.Range("I104").Value = 2800
$ws.Range("J104").Value = 4457.25
$ws.Range("K104").Value = 8400
$ws.Range("L104").Value = 13371.75
$ws.Range("M104").Value = -5779
$ws.Range("N104").Value = -18613.75
$ws.Range("H122").Value = 4291.143
$ws.Range("J122").Value = 7633.533
$ws.Range("L122").Value = 68701.79700000001
$ws.Range("N122").Value = -73601.79700000001
$ws.Range("H132").Value = 2634.9707
$ws.Range("I132").Value = 2681.75
$ws.Range("J132").Value = 2620.577
$ws.Range("K132").Value = 24135.75
$ws.Range("L132").Value = 23585.193
$ws.Range("M132").Value = -21605.75
$ws.Range("N132").Value = -28645.193
$ws.Range("H135").Value = 802.15
$ws.Range("J135").Value = 1243.2858
$ws.Range("L135").Value = 11189.5722
$ws.Range("N135").Value = -16259.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1843.3
$ws.Range("I122").Value = 1738.8334
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5216.5002
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2766.5002
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130
$ws.Range("H132").Value = 3259.7632
$ws.Range("I132").Value = 2870.6553
$ws.Range("J132").Value = 4513.5557
$ws.Range("K132").Value = 8611.965899999999
$ws.Range("L132").Value = 13540.6671
$ws.Range("M132").Value = -6081.965899999999
$ws.Range("N132").Value = -18600.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3769.2104
$ws.Range("I7").Value = 2545.889
$ws.Range("J7").Value = 4870.2
$ws.Range("K7").Value = 2545.889
$ws.Range("L7").Value = 4870.2
$ws.Range("M7").Value = -2433.889
$ws.Range("N7").Value = -5094.2
$ws.Range("H75").Value = 55782
$ws.Range("J75").Value = 55782
$ws.Range("L75").Value = 55782
$ws.Range("N75").Value = -57654
$ws.Range("H78").Value = 55782
$ws.Range("J78").Value = 55782
$ws.Range("L78").Value = 167346
$ws.Range("N78").Value = -176706
$ws.Range("H126").Value = 3769.2104
$ws.Range("I126").Value = 2545.889
$ws.Range("J126").Value = 4870.2
$ws.Range("K126").Value = 7637.667
$ws.Range("L126").Value = 14610.6
$ws.Range("M126").Value = -5167.667
$ws.Range("N126").Value = -19550.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 70049
$ws.Range("J42").Value = 70049
$ws.Range("L42").Value = 70049
$ws.Range("N42").Value = -70805
$ws.Range("H132").Value = 3088164.2
$ws.Range("I132").Value = 1796.375
$ws.Range("K132").Value = 5389.125
$ws.Range("M132").Value = -2859.125
$ws.Range("H136").Value = 2438.1177
$ws.Range("I136").Value = 1957.1786
$ws.Range("J136").Value = 4682.5
$ws.Range("K136").Value = 5871.5358
$ws.Range("L136").Value = 14047.5
$ws.Range("M136").Value = -3321.5358
$ws.Range("N136").Value = -19147.5
